$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) The "Date 3" placeholder (an auto-updating datetimeFigureOut field) on
#    the slide master and on every slide layout shows a cached date string.
#    Bump it from 26.10.2023 to 01.11.2023 everywhere it appears.
# ---------------------------------------------------------------------------
function Set-DatePlaceholderText($shapes, [string]$newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePlaceholder = $false
        if ($shp.Type -eq 14) {
            try {
                if ($shp.PlaceholderFormat.Type -eq 16) {
                    $isDatePlaceholder = $true
                }
            } catch {
            }
        }
        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

$newDate = "01.11.2023"

# Slide master
Set-DatePlaceholderText $p.SlideMaster.Shapes $newDate

# Every slide layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Set-DatePlaceholderText $layouts.Item($li).Shapes $newDate
}

# ---------------------------------------------------------------------------
# 2) Nudge "TextBox 6" on slide 3 to the left (763570 EMU -> 730576 EMU).
#    PowerPoint's COM object model works in points (1 pt = 12700 EMU), so
#    730576 / 12700 = 57.5257 pt.
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$textBox6 = $slide3.Shapes.Item(3)
$textBox6.Left = 57.5257
